# Auto update ESO load data
# Updates the daily load-forecast table on Sheet1 (rows 2-7, columns A:AA)
# with the latest rolled-forward values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(2026,1,31,4758,4540,4440,4405,4426,4488,4612,4857,5365,5795,5920,5848,5649,5449,5278,5257,5409,5790,6007,5891,5699,5486,5373,5174),
    @(2026,2,1,4876,4744,4700,4659,4614,4691,4801,5106,5536,6003,6212,6162,6066,5890,5750,5764,5967,6366,6623,6559,6388,6130,5949,5654),
    @(2026,2,2,5308,5147,5072,5073,5152,5314,5817,6408,6924,7094,6921,6734,6522,6442,6373,6406,6439,6781,7223,7198,7019,6769,6542,6183),
    @(2026,2,3,5769,5504,5430,5424,5463,5616,6111,6676,7151,7270,7123,6969,6761,6648,6527,6533,6625,6882,7196,7183,6983,6725,6486,6140),
    @(2026,2,4,5693,5437,5366,5361,5396,5535,5986,6500,6933,7041,6907,6767,6576,6473,6363,6369,6452,6688,6975,6962,6759,6498,6255,5905),
    @(2026,2,5,5370,5117,5046,5040,5076,5214,5663,6174,6605,6712,6597,6476,6312,6223,6128,6133,6205,6409,6657,6644,6437,6172,5925,5569)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
